$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nädal 9")

# Row 4
$ws.Range("F4").Value = 80
$ws.Range("G4").Value = "video"

# Row 5
$ws.Range("D5").Value = 0.625
$ws.Range("F5").Value = 30
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("G5").Value = "video"
$ws.Range("J5").Value = 15

# Row 6
$ws.Range("B6").Value = 43919
$ws.Range("C6").Value = 0.57291666666666663
$ws.Range("D6").Value = 0.60416666666666663
$ws.Range("E6").Value = 15
$ws.Range("F6").Value = 30
$ws.Range("G4").Copy()
$ws.Range("G6").PasteSpecial(-4122)
$ws.Range("G6").Value = "video"
$ws.Range("H6").Value = "V40"
$ws.Range("J6").Value = 46

# Row 7
$ws.Range("C7").Value = 0.72916666666666663
$ws.Range("D7").Value = 0.875
$ws.Range("E7").Value = 30
$ws.Range("F7").Value = 180
$ws.Range("G4").Copy()
$ws.Range("G7").PasteSpecial(-4122)
$ws.Range("G7").Value = "video"
$ws.Range("H7").Value = "V40, 41"
$ws.Range("I7").Value = "x"

# Row 8
$ws.Range("C8").Value = 0.875
$ws.Range("G4").Copy()
$ws.Range("G8").PasteSpecial(-4122)
$ws.Range("G8").Value = "video"

# Row 9
$ws.Range("G4").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Value = "video"

# Row heights / thick bottom borders for rows 4-8 (match formatting from diff)
$ws.Range("A4:J8").RowHeight = 15.75

# Selection
$ws.Range("D8").Select()
